# Weekly fruit/vegetable price update: a new daily record is inserted
# for "Macroferia Regional de Talca - Betarraga" (Hortaliza) ahead of the
# existing row 148, pushing all subsequent records (148-182) down by one
# row (to 149-183) and growing the used range from A1:R182 to A1:R183.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 148; Excel shifts rows 148:182 -> 149:183
# and carries formatting (e.g. the date style on column D) down with them.
$ws.Rows.Item(148).Insert()

# Populate the newly inserted row 148 with the new weekly observation.
$ws.Cells.Item(148, 1).Value = 5
$ws.Cells.Item(148, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(148, 3).Value = "Maule"
$ws.Cells.Item(148, 4).Value = 44476
$ws.Cells.Item(148, 5).Value = 7
$ws.Cells.Item(148, 6).Value = 100114014
$ws.Cells.Item(148, 7).Value = "Betarraga"
$ws.Cells.Item(148, 8).Value = "Sin especificar"
$ws.Cells.Item(148, 9).Value = "Primera"
$ws.Cells.Item(148, 10).Value = 4000
$ws.Cells.Item(148, 11).Value = 650
$ws.Cells.Item(148, 12).Value = 650
$ws.Cells.Item(148, 13).Value = 650
$ws.Cells.Item(148, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(148, 15).Value = "Región del Maule"
$ws.Cells.Item(148, 16).Value = 130
$ws.Cells.Item(148, 17).Value = 5
$ws.Cells.Item(148, 18).Value = "Hortaliza"
